$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.293.66"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.827.55"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'363.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Value = "'111.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.50%  "
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'40.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "'19.88"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'7.82"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "3.283.57"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "2.865.39"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "'0.928"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").Value = "52.158.03"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "'3.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "'13.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").Value = "0.0₃0998"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "'272.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "'70.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'26.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'10.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'0.0475"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").Value = "'52.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("D33").Value = "'35.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "'5.54"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.66%  "
$ws.Range("D36").Value = "'0.0851"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'3.26"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'2.03"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("D40").Value = "'18.31"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'2.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Value = "'125.47"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").Value = "'22.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").Value = "2.078.37"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'2.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "'5.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D50").Value = "'0.961"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").Value = "'9.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.61%  "
